$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 update (unmapped) - F,G,H,I literal values change, J becomes formula, K:Q shared formula, R new avg
$ws.Range("F2").Value = 2230250
$ws.Range("G2").Value = 2464348
$ws.Range("H2").Value = 1355386
$ws.Range("I2").Value = 3092229
$ws.Range("J2").Formula = "=J15-J4"
$ws.Range("K2:Q2").Formula = "=K15-K4"
$ws.Range("R2").Formula = "=AVERAGE(F2:Q2)"
$ws.Range("AH2").Formula = "=AG13"

# Row 3 update (unclassified) - F3=0, clear G3:Q3
$ws.Range("F3").Value = 0
$ws.Range("G3:Q3").ClearContents()

# Row 4 update (up to species)
$ws.Range("F4").Value = 16044
$ws.Range("G4").Formula = "=36473+7213+2401+14763+406153"
$ws.Range("H4").Formula = "=15111"
$ws.Range("I4").Formula = "=13160+29812"
$ws.Range("J4").Formula = "=2470+66877"
$ws.Range("K4").Formula = "=30600+667+22"
$ws.Range("L4").Formula = "=401383+12099+1629+6+45027+10542+2"
$ws.Range("M4").Formula = "=25712+1104+28894"
$ws.Range("N4").Formula = "=5278+52296+9+697+5776+160304"
$ws.Range("O4").Formula = "=31820+10721"
$ws.Range("P4").Formula = "=16000+19+7109+150652"
$ws.Range("Q4").Formula = "=3+4549+39780+6+1723+7062+146996"
$ws.Range("R4").Formula = "=AVERAGE(F4:Q4)"

# Row 11 (NEW): total annotated = SUM(row4:row10) per column
$ws.Range("A11").Value = "total annotated"
$ws.Range("F11").Formula = "=SUM(F4:F10)"
$ws.Range("G11:R11").Formula = "=SUM(G4:G10)"
$ws.Range("U11").Formula = "=SUM(U4:U10)"
$ws.Range("V11:AG11").Formula = "=SUM(V4:V10)"

# Row 13 (NEW, content = old row12's formulas): "total"
$ws.Range("A13").Value = "total"
$ws.Range("F13").Formula = "=SUM(F2:F4)"
$ws.Range("G13:R13").Formula = "=SUM(G2:G4)"
$ws.Range("U13").Formula = "=SUM(U2:U10)"
$ws.Range("V13:AF13").Formula = "=SUM(V2:V10)"
$ws.Range("AG13").Formula = "=SUM(AG2:AG10)"

# Row 14 (was row13 "raw reads") - literal values unchanged from old row13
$ws.Range("A14").Value = "raw reads"
$ws.Range("F14").Value = 8985176
$ws.Range("G14").Value = 11725404
$ws.Range("H14").Value = 5481988
$ws.Range("I14").Value = 12540804
$ws.Range("J14").Value = 8086784
$ws.Range("K14").Value = 6317488
$ws.Range("L14").Value = 12350020
$ws.Range("M14").Value = 6883748
$ws.Range("N14").Value = 10707296
$ws.Range("O14").Value = 10089840
$ws.Range("P14").Value = 16571364
$ws.Range("Q14").Value = 9383212

# Row 15 (was row14 "wc") - now referencing row14
$ws.Range("A15").Value = "wc"
$ws.Range("F15").Formula = "=F14/4"
$ws.Range("G15:Q15").Formula = "=G14/4"
$ws.Range("R15").Formula = "=AVERAGE(F15:Q15)"
$ws.Range("U15").Formula = "=F15"
$ws.Range("V15:AF15").Formula = "=G15"
$ws.Range("AG15").Formula = "=AVERAGE(U15:AF15)"

# Row 12 (discarded) - depends on row13(total) and row15(wc)
$ws.Range("A12").Value = "discarded"
$ws.Range("F12").Formula = "=F13-F15"
$ws.Range("G12:Q12").Formula = "=G13-G15"
$ws.Range("U12").Formula = "=U15-U13"
$ws.Range("V12:AG12").Formula = "=V15-V13"

# Row 1 new header R1
$ws.Range("R1").Value = "humann2 avg"

# Rows 18-28: summary table for the chart
$ws.Range("B18").Value = "humann2 avg"
$ws.Range("C18").Value = "mpro avg"

$ws.Range("A19").Value = "discarded"
$ws.Range("B19").Formula = "=0"
$ws.Range("C19").Formula = "=AG12"

$ws.Range("A20").Value = "unmapped"
$ws.Range("B20:B28").Formula = "=R2"
$ws.Range("C20:C28").Formula = "=AG2"

$ws.Range("A21").Value = "unclassified"
$ws.Range("A22").Value = "up to species"
$ws.Range("A23").Value = "up to genus"
$ws.Range("A24").Value = "up to family"
$ws.Range("A25").Value = "up to order"
$ws.Range("A26").Value = "up to class"
$ws.Range("A27").Value = "up to phylum"
$ws.Range("A28").Value = "up to kingdom"

Write-Host "all rows done"

$ws.Calculate()
Write-Host "=== Rows 18-28 Verification ==="
for ($r = 18; $r -le 28; $r++) {
    $a = $ws.Range("A$r").Value2
    $b = $ws.Range("B$r").Value2
    $c = $ws.Range("C$r").Value2
    Write-Host "Row$r A=$a B=$b C=$c"
}
